# Update the worksheet date and all 25 division-problem answers.
# The table's shape (20 rows x 5 columns = 100 cells) does not change;
# only the text content of the date paragraph and the cells changes.
# Two of the new values collide with an existing value elsewhere in the
# document (row 5, col 3 currently reads "62÷3=20, 2", which is also the
# text that row 5, col 1 must become), so those two cells are addressed
# directly via the table/row/cell object model instead of Find/Replace
# to avoid any ambiguity. All other changes use Find & Replace, which is
# safe here because every "find" string is unique in the document.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Header date line.
Replace-Text "2023-08-31 Thursday" "2023-09-01 Friday"

# Row 1 (unambiguous simple replacements).
Replace-Text "33÷4=8, 1" "82÷4=20, 2"
Replace-Text "18÷6=3, 0" "88÷4=22, 0"
Replace-Text "60÷2=30, 0" "50÷8=6, 2"
Replace-Text "73÷5=14, 3" "79÷8=9, 7"
Replace-Text "87÷2=43, 1" "90÷5=18, 0"

# Row 5 - handled directly by cell position because the new text for
# column 1 collides with the pre-existing text in column 3.
$table = $d.Tables.Item(1)
$row5 = $table.Rows.Item(5)
$row5.Cells.Item(1).Range.Text = "62÷3=20, 2"   # was 78÷3=26, 0
$row5.Cells.Item(2).Range.Text = "54÷5=10, 4"   # was 88÷7=12, 4
$row5.Cells.Item(3).Range.Text = "86÷2=43, 0"   # was 62÷3=20, 2
$row5.Cells.Item(4).Range.Text = "84÷4=21, 0"   # was 57÷4=14, 1
$row5.Cells.Item(5).Range.Text = "48÷3=16, 0"   # was 55÷4=13, 3

# Row 9.
Replace-Text "27÷6=4, 3" "62÷4=15, 2"
Replace-Text "92÷8=11, 4" "52÷4=13, 0"
Replace-Text "71÷8=8, 7" "34÷2=17, 0"
Replace-Text "57÷3=19, 0" "74÷9=8, 2"
Replace-Text "65÷6=10, 5" "22÷9=2, 4"

# Row 13.
Replace-Text "35÷8=4, 3" "53÷3=17, 2"
Replace-Text "47÷2=23, 1" "35÷2=17, 1"
Replace-Text "18÷5=3, 3" "62÷7=8, 6"
Replace-Text "39÷8=4, 7" "97÷4=24, 1"
Replace-Text "35÷4=8, 3" "79÷3=26, 1"

# Row 17.
Replace-Text "10÷7=1, 3" "24÷8=3, 0"
Replace-Text "69÷4=17, 1" "66÷7=9, 3"
Replace-Text "53÷4=13, 1" "26÷6=4, 2"
Replace-Text "40÷5=8, 0" "48÷9=5, 3"
Replace-Text "86÷3=28, 2" "38÷8=4, 6"
